# Weekly Fruta/Hortaliza update: insert two new daily-price records
# (Cilantro, Terminal Hortofruticola Agro Chillan) for the newest reporting
# date, pushing the existing rows 134-242 down to 136-244.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 134-135; everything from the old row 134 down
# shifts to rows 136-244 (dimension grows from A1:R242 to A1:R244).
$ws.Range("A134:A135").EntireRow.Insert()

# ---- New row 134: Cilantro, Primera ----
$ws.Cells.Item(134, 1).Value = 7
$ws.Cells.Item(134, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(134, 3).Value = "Ñuble"
$ws.Cells.Item(134, 4).Value = 45072
$ws.Cells.Item(134, 5).Value = 16
$ws.Cells.Item(134, 6).Value = 100112040
$ws.Cells.Item(134, 7).Value = "Cilantro"
$ws.Cells.Item(134, 8).Value = "Sin especificar"
$ws.Cells.Item(134, 9).Value = "Primera"
$ws.Cells.Item(134, 10).Value = 100
$ws.Cells.Item(134, 11).Value = 1200
$ws.Cells.Item(134, 12).Value = 1200
$ws.Cells.Item(134, 13).Value = 1200
$ws.Cells.Item(134, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(134, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(134, 16).Value = 1200
$ws.Cells.Item(134, 17).Value = 1
$ws.Cells.Item(134, 18).Value = "Hortaliza"

# ---- New row 135: Cilantro, Segunda ----
$ws.Cells.Item(135, 1).Value = 7
$ws.Cells.Item(135, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(135, 3).Value = "Ñuble"
$ws.Cells.Item(135, 4).Value = 45072
$ws.Cells.Item(135, 5).Value = 16
$ws.Cells.Item(135, 6).Value = 100112040
$ws.Cells.Item(135, 7).Value = "Cilantro"
$ws.Cells.Item(135, 8).Value = "Sin especificar"
$ws.Cells.Item(135, 9).Value = "Segunda"
$ws.Cells.Item(135, 10).Value = 100
$ws.Cells.Item(135, 11).Value = 1000
$ws.Cells.Item(135, 12).Value = 1000
$ws.Cells.Item(135, 13).Value = 1000
$ws.Cells.Item(135, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(135, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(135, 16).Value = 1000
$ws.Cells.Item(135, 17).Value = 1
$ws.Cells.Item(135, 18).Value = "Hortaliza"
